$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.122.12"
$ws.Range("E2").Value = "  -1.19%  "
$ws.Range("D3").Value = "3.204.80"
$ws.Range("E3").Value = "  -1.64%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.73%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.597"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.41%  "
$ws.Range("E9").Value = "  -4.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.58"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.409"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.52%  "
$ws.Range("D12").Value = "3.783.54"
$ws.Range("E12").Value = "  -1.03%  "
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.63%  "
$ws.Range("D15").Value = "67.271.41"
$ws.Range("E15").Value = "  -0.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000167"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.94%  "
$ws.Range("D17").Value = "3.219.70"
$ws.Range("E17").Value = "  -1.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.70"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.33"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "390.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.509"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000116"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.184"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.04%  "
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.94"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.46"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.91"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.48%  "
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.23"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "160.76"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.46"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.86"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.42%  "
$ws.Range("B38").Value = "Mantle"
$ws.Range("C38").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.802"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.90%  "
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "26.14"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.50"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "40.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.50%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.43"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.66%  "
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0676"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.85%  "
$ws.Range("D45").Value = "2.584.80"
$ws.Range("E45").Value = "  -1.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "24.58"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "330.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0274"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.21"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.100"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.81%  "
$ws.Range("B51").Value = "FirstDigitalUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.05%  "
